# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.462.24"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.567.69"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.17"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.790.85"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "1.554.12"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.69"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "27.466.00"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").Value = "0.0₃0690"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "1.377.02"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.958"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.547"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.828"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.31"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "1.703.30"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.27"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "0.0₇0999"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0958"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -0.52%  "
